# Apply room/enrollment reshuffles and add new ENS207 course row.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-Cell($row, $col, $value) {
    $ws.Cells.Item($row, $col).Value = $value
}

# Simple room-name swaps (and a few enrollment/capacity tweaks)
Set-Cell 6   2 "B F2.5"

Set-Cell 15  2 "A F2.16 - Architecture Studio"
Set-Cell 15  5 20

Set-Cell 37  2 "A B.13 - Class/PSY Lab"
Set-Cell 37  5 12

Set-Cell 39  2 "A F3.8 - Big Architecture Studio"

Set-Cell 60  2 "B F2.15 - Amphitheater II"

Set-Cell 67  2 "A F1.25"

Set-Cell 84  2 "A F1.3 - Computer Lab"
Set-Cell 84  5 25

Set-Cell 98  2 "RC1.4 - Computer Laboratory"

Set-Cell 105 2 "RC1.3 - GSM and Network Laboratories"

Set-Cell 106 2 "A F1.18 - Computer Lab"

Set-Cell 124 2 "B F1.2 - Class/ECON Lab"

# Row 131: course renamed from ENS207-6.1 to ENS209 with new room/time/counts
Set-Cell 131 1 "ENS209"
Set-Cell 131 2 "B F1.16"
Set-Cell 131 3 "Tue. 17:00-19:50"
Set-Cell 131 4 39
Set-Cell 131 5 40

Set-Cell 134 2 "A F1.17"

Set-Cell 136 2 "B F1.10 Class/ART Studio"

Set-Cell 148 2 "A F3.7 - Small Architecture Studio & A F3.10 - Architecture Classroom"

Set-Cell 149 2 "B F1.23 - Amphitheater I"

Set-Cell 155 2 "A F2.8 - Drawing Studio"

Set-Cell 163 2 "B F1.10 Class/ART Studio"

Set-Cell 165 2 "B F1.25 Computer Lab"

Set-Cell 166 2 "A F2.8 - Drawing Studio & A F2.16 - Architecture Studio"

Set-Cell 172 2 "B F1.35 FBA Conference Room"

Set-Cell 175 2 "RC1.4 - Computer Laboratory"

Set-Cell 178 2 "A B.2 - EE Lab"

Set-Cell 183 2 "A F1.4 - Class/Laboratory"

Set-Cell 184 2 "B F1.10 Class/ART Studio"

Set-Cell 187 2 "B F1.35 FBA Conference Room"

Set-Cell 188 2 "B F2.17"

Set-Cell 192 2 "RC1.3 - GSM and Network Laboratories"

Set-Cell 196 2 "B F2.15 - Amphitheater II"

Set-Cell 197 2 "B F1.23 - Amphitheater I"

Set-Cell 199 2 "B F1.23 - Amphitheater I"

Set-Cell 205 2 "A F1.4 - Class/Laboratory"

Set-Cell 212 2 "A F1.10"

Set-Cell 219 2 "B F1.10 Class/ART Studio"

Set-Cell 224 2 "B F2.5"

Set-Cell 233 2 "B F1.16"

Set-Cell 264 2 "B F1.10 Class/ART Studio"

Set-Cell 265 2 "B F2.2"

Set-Cell 284 2 "B F1.2 - Class/ECON Lab"

Set-Cell 304 2 "B F2.27 Creative Writing and Translation Studio"
Set-Cell 304 5 18

Set-Cell 305 2 "B F2.27 Creative Writing and Translation Studio"
Set-Cell 305 5 18

Set-Cell 306 2 "A F3.8 - Big Architecture Studio"
Set-Cell 306 5 25

Set-Cell 307 2 "RC1.3 - GSM and Network Laboratories"
Set-Cell 307 5 20

Set-Cell 314 2 "A B.2 - EE Lab"

Set-Cell 319 2 "RC1.4 - Computer Laboratory"

Set-Cell 325 2 "A F2.8 - Drawing Studio"
Set-Cell 325 5 25

Set-Cell 326 2 "A B.13 - Class/PSY Lab"
Set-Cell 326 5 12

Set-Cell 327 2 "RC1.4 - Computer Laboratory"
Set-Cell 327 5 20

Set-Cell 329 2 "B F2.27 Creative Writing and Translation Studio"
Set-Cell 329 5 18

Set-Cell 332 2 "RC1.3 - GSM and Network Laboratories"
Set-Cell 332 5 20

# New row 334: ENS207 course
Set-Cell 334 1 "ENS207"
Set-Cell 334 2 "B F2.16"
Set-Cell 334 3 "Fri. 09:00-09:50"
Set-Cell 334 4 71
Set-Cell 334 5 80
Set-Cell 334 6 "Assigned"
